# Update Work Week and Social Spending
# - Replace the GDP-per-Capita-era data values in column E (Data sheet)
#   with the new "Work Week and Social Spending" series for 1820-2008.
# - Append 8 new rows (years 2009-2016) with the new series values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# (row, new value) pairs for existing rows whose Data value changed
$changedData = @(
    @(2, "641"),
    @(52, "641"),
    @(95, "870"),
    @(132, "802"),
    @(133, "816"),
    @(134, "835"),
    @(135, "877"),
    @(136, "886"),
    @(137, "894"),
    @(138, "925"),
    @(139, "915"),
    @(140, "956"),
    @(141, "971"),
    @(142, "980"),
    @(143, "990"),
    @(144, "998"),
    @(145, "1006"),
    @(146, "1012"),
    @(147, "1019"),
    @(148, "1071"),
    @(149, "1034"),
    @(150, "1023"),
    @(151, "1049"),
    @(152, "1055"),
    @(153, "1022"),
    @(154, "1034"),
    @(155, "1009"),
    @(156, "1052"),
    @(157, "1047"),
    @(158, "1071"),
    @(159, "1081"),
    @(160, "1105"),
    @(161, "1106"),
    @(162, "1057"),
    @(163, "1119"),
    @(164, "1135"),
    @(165, "1076"),
    @(166, "1152"),
    @(167, "1194"),
    @(168, "1218"),
    @(169, "1208"),
    @(170, "1269"),
    @(171, "1290"),
    @(172, "1315"),
    @(173, "1360.78354501199"),
    @(174, "1374.66118498648"),
    @(175, "1385.24106157196"),
    @(176, "1456.02029656526"),
    @(177, "1468.07894532892"),
    @(178, "1511.73367300972"),
    @(179, "1558.39878541556"),
    @(180, "1570.66349301421"),
    @(181, "1609.18027622043"),
    @(182, "1676.7456451079"),
    @(183, "1738.31417941039"),
    @(184, "1709.77869937907"),
    @(185, "1747.77660138511"),
    @(186, "1802.14384785063"),
    @(187, "1836.30571517462"),
    @(188, "1884.18716679371"),
    @(189, "1927.38650311044"),
    @(190, "2021.15408379838")
)

# (row, year, value) triples for brand-new rows (2009-2016)
$newRowsData = @(
    @(191, 2009, "2085.24361411454"),
    @(192, 2010, "2152.27458744838"),
    @(193, 2011, "2189"),
    @(194, 2012, "2258"),
    @(195, 2013, "2316"),
    @(196, 2014, "2421"),
    @(197, 2015, "2455"),
    @(198, 2016, "2435")
)

foreach ($pair in $changedData) {
    $row = $pair[0]
    $value = $pair[1]
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

foreach ($triple in $newRowsData) {
    $row = $triple[0]
    $year = $triple[1]
    $value = $triple[2]

    $ws.Cells.Item($row, 1).Value = 524
    $ws.Cells.Item($row, 2).Value = "Nepal"
    $ws.Cells.Item($row, 3).Value = "GDP per Capita"
    $ws.Cells.Item($row, 4).Value = $year

    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}
